$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.189.57"
$ws.Range("E2").Value = "  -3.77%  "
$ws.Range("D3").Value = "3.133.46"
$ws.Range("E3").Value = "  -5.45%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.23%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.132.30"
$ws.Range("E8").Value = "  -5.49%  "
$ws.Range("E9").Value = "  -5.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.50%  "
$ws.Range("E11").Value = "  -8.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.90%  "
$ws.Range("D13").Value = "3.665.24"
$ws.Range("E13").Value = "  -5.52%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("D16").Value = "3.128.10"
$ws.Range("E16").Value = "  -5.49%  "
$ws.Range("D17").Value = "58.083.52"
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("E18").Value = "  -8.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.08%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.18%  "
$ws.Range("E25").Value = "  -6.32%  "
$ws.Range("D26").Value = "3.252.85"
$ws.Range("E26").Value = "  -5.59%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0955"
$ws.Range("E27").Value = "  -6.99%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.167"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -9.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.48%  "
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0691"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "3.162.38"
$ws.Range("E41").Value = "  -5.40%  "
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.696"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.63%  "
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.40%  "
$ws.Range("D49").Value = "2.259.70"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.24%  "

Write-Output "done"